$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H92").Value = 912.9286
$ws.Range("I92").Value = 648
$ws.Range("K92").Value = 648
$ws.Range("M92").Value = 600
$ws.Range("H96").Value = 277.07144
$ws.Range("I96").Value = 236.84616
$ws.Range("J96").Value = 800
$ws.Range("K96").Value = 710.5384799999999
$ws.Range("L96").Value = 2400
$ws.Range("M96").Value = 662.4615200000001
$ws.Range("N96").Value = -5146
$ws.Range("H100").Value = 2034.7778
$ws.Range("I100").Value = 1058.2858
$ws.Range("J100").Value = 5452.5
$ws.Range("K100").Value = 1058.2858
$ws.Range("L100").Value = 5452.5
$ws.Range("M100").Value = -517.2858000000001
$ws.Range("N100").Value = -6534.5
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 1258800.8
$ws.Range("I132").Value = 2322.6128
$ws.Range("K132").Value = 6967.8384
$ws.Range("M132").Value = -4437.8384
$ws.Range("H137").Value = 3706812.5
$ws.Range("I137").Value = 6668949.5
$ws.Range("J137").Value = 4141.3335
$ws.Range("K137").Value = 20006848.5
$ws.Range("L137").Value = 12424.0005
$ws.Range("M137").Value = -20004298.5
$ws.Range("N137").Value = -17524.0005
$ws.Range("H141").Value = 649.2759
$ws.Range("I141").Value = 618.8929000000001
$ws.Range("J141").Value = 1500
$ws.Range("K141").Value = 1856.6787
$ws.Range("L141").Value = 4500
$ws.Range("M141").Value = 3323.3213
$ws.Range("N141").Value = -14860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 903.97
$ws.Range("I32").Value = 788.9036
$ws.Range("K32").Value = 788.9036
$ws.Range("M32").Value = -501.9036
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H63").Value = 15000
$ws.Range("I63").Value = 10000
$ws.Range("K63").Value = 10000
$ws.Range("M63").Value = -9314
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H66").Value = 15000
$ws.Range("I66").Value = 10000
$ws.Range("K66").Value = 50000
$ws.Range("M66").Value = -46568
$ws.Range("H102").Value = 7525774.5
$ws.Range("I102").Value = 8936295
$ws.Range("K102").Value = 8936295
$ws.Range("M102").Value = -8934673
$ws.Range("H112").Value = 27199.8
$ws.Range("J112").Value = 27199.8
$ws.Range("L112").Value = 27199.8
$ws.Range("N112").Value = -30153.8
$ws.Range("H125").Value = 54425.625
$ws.Range("J125").Value = 54425.625
$ws.Range("L125").Value = 54425.625
$ws.Range("N125").Value = -64265.625
$ws.Range("H132").Value = 78523.7
$ws.Range("I132").Value = 50830.55
$ws.Range("K132").Value = 152491.65
$ws.Range("M132").Value = -149961.65
$ws.Range("H133").Value = 32156.6
$ws.Range("J133").Value = 32156.6
$ws.Range("L133").Value = 32156.6
$ws.Range("N133").Value = -37216.6
$ws.Range("H135").Value = 49736.844
$ws.Range("J135").Value = 49736.844
$ws.Range("L135").Value = 49736.844
$ws.Range("N135").Value = -59876.844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 40181
$ws.Range("J62").Value = 40181
$ws.Range("L62").Value = 40181
$ws.Range("N62").Value = -41553
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H65").Value = 40181
$ws.Range("J65").Value = 40181
$ws.Range("L65").Value = 120543
$ws.Range("N65").Value = -127407
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H94").Value = 718.9
$ws.Range("I94").Value = 648.5833
$ws.Range("K94").Value = 648.5833
$ws.Range("M94").Value = -197.5833
$ws.Range("H99").Value = 1259.3182
$ws.Range("I99").Value = 1195.4736
$ws.Range("K99").Value = 1195.4736
$ws.Range("M99").Value = 302.5264
$ws.Range("H110").Value = 30000
$ws.Range("J110").Value = 30000
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180
$ws.Range("H111").Value = 38600.668
$ws.Range("J111").Value = 38600.668
$ws.Range("L111").Value = 38600.668
$ws.Range("N111").Value = -46780.668
$ws.Range("H134").Value = 3470.0625
$ws.Range("I134").Value = 2300.1428
$ws.Range("J134").Value = 4380
$ws.Range("K134").Value = 6900.428400000001
$ws.Range("L134").Value = 13140
$ws.Range("M134").Value = -4365.428400000001
$ws.Range("N134").Value = -18210

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1225.25
$ws.Range("J16").Value = 995
$ws.Range("L16").Value = 995
$ws.Range("N16").Value = -1569
$ws.Range("H31").Value = 1917.1187
$ws.Range("I31").Value = 1258.5454
$ws.Range("J31").Value = 3848.9333
$ws.Range("K31").Value = 1258.5454
$ws.Range("L31").Value = 3848.9333
$ws.Range("M31").Value = -963.5454
$ws.Range("N31").Value = -4438.933300000001
$ws.Range("H34").Value = 1917.1187
$ws.Range("I34").Value = 1258.5454
$ws.Range("J34").Value = 3848.9333
$ws.Range("K34").Value = 1258.5454
$ws.Range("L34").Value = 3848.9333
$ws.Range("M34").Value = -1056.5454
$ws.Range("N34").Value = -4252.933300000001
$ws.Range("H113").Value = 1225.25
$ws.Range("J113").Value = 995
$ws.Range("L113").Value = 995
$ws.Range("N113").Value = -5335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 24390310
$ws.Range("I12").Value = 52631664
$ws.Range("J12").Value = 51.545456
$ws.Range("K12").Value = 157894992
$ws.Range("L12").Value = 154.636368
$ws.Range("M12").Value = -157894819
$ws.Range("N12").Value = -500.636368
$ws.Range("H80").Value = 1275.5238
$ws.Range("J80").Value = 1464.3334
$ws.Range("L80").Value = 4393.0002
$ws.Range("N80").Value = -6265.0002
$ws.Range("H83").Value = 1275.5238
$ws.Range("J83").Value = 1464.3334
$ws.Range("L83").Value = 13179.0006
$ws.Range("N83").Value = -22539.0006
$ws.Range("H92").Value = 998.9091
$ws.Range("I92").Value = 998.9091
$ws.Range("K92").Value = 2996.7273
$ws.Range("M92").Value = -1748.7273
$ws.Range("H122").Value = 698.03845
$ws.Range("J122").Value = 1282.909
$ws.Range("L122").Value = 11546.181
$ws.Range("N122").Value = -16446.181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 28561.75
$ws.Range("J92").Value = 28082.334
$ws.Range("L92").Value = 28082.334
$ws.Range("N92").Value = -31826.334
$ws.Range("H123").Value = 22721.941
$ws.Range("J123").Value = 22721.941
$ws.Range("L123").Value = 22721.941
$ws.Range("N123").Value = -27621.941

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2048.75
$ws.Range("I7").Value = 2036.3636
$ws.Range("J7").Value = 2076
$ws.Range("K7").Value = 2036.3636
$ws.Range("L7").Value = 2076
$ws.Range("M7").Value = -1924.3636
$ws.Range("N7").Value = -2300
$ws.Range("H40").Value = 2639.8667
$ws.Range("I40").Value = 2724.75
$ws.Range("K40").Value = 2724.75
$ws.Range("M40").Value = -2588.75
$ws.Range("H93").Value = 313.5
$ws.Range("H111").Value = 60000
$ws.Range("J111").Value = 60000
$ws.Range("L111").Value = 60000
$ws.Range("N111").Value = -68180
$ws.Range("H126").Value = 2048.75
$ws.Range("I126").Value = 2036.3636
$ws.Range("J126").Value = 2076
$ws.Range("K126").Value = 6109.0908
$ws.Range("L126").Value = 6228
$ws.Range("M126").Value = -3639.0908
$ws.Range("N126").Value = -11168
$ws.Range("H136").Value = 155718
$ws.Range("I136").Value = 186175.67
$ws.Range("J136").Value = 132874.75
$ws.Range("K136").Value = 558527.01
$ws.Range("L136").Value = 398624.25
$ws.Range("M136").Value = -555977.01
$ws.Range("N136").Value = -403724.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1700
$ws.Range("I96").Value = 1533.8334
$ws.Range("J96").Value = 1842.4286
$ws.Range("K96").Value = 1533.8334
$ws.Range("L96").Value = 1842.4286
$ws.Range("M96").Value = -160.8334
$ws.Range("N96").Value = -4588.4286
$ws.Range("H100").Value = 67716.53
$ws.Range("I100").Value = 100496
$ws.Range("K100").Value = 200992
$ws.Range("M100").Value = -200451
$ws.Range("H126").Value = 1490.1111
$ws.Range("I126").Value = 1127.75
$ws.Range("J126").Value = 1780
$ws.Range("K126").Value = 3383.25
$ws.Range("L126").Value = 5340
$ws.Range("M126").Value = -913.25
$ws.Range("N126").Value = -10280
